$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new paragraph right after the blank line that follows
#    "Currently Supported Platforms: ..." containing the new "Download
#    Respective Build..." text, followed by a blank paragraph, before
#    "Instructions for Windows Users:". The new paragraph also carries the
#    "_GoBack" bookmark that used to sit at the end of the document.
# ---------------------------------------------------------------------------

$introBlank = $d.Paragraphs.Item(3)
$rng = $introBlank.Range
$insertionPoint = $d.Range($rng.End, $rng.End)
# Trailing placeholder "X" keeps the collapsed bookmark range away from the
# paragraph-final boundary (a spot where this engine mis-resolves zero
# length ranges); it is stripped again once the bookmark is anchored.
$insertionPoint.Text = "Download Respective Build for either Windows or Mac OS from Builds FolderX`r"

$newPara = $d.Paragraphs.Item(4)
$newPara.Range.InsertParagraphAfter()

$foundRange = $d.Content
$foundRange.Find.Execute("Download Respective Build for either Windows or Mac OS from Builds FolderX", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkSpot = $d.Range($foundRange.End - 1, $foundRange.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null

$cleanupRange = $d.Content
$cleanupRange.Find.Execute("FolderX", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$trailingChar = $d.Range($cleanupRange.End - 1, $cleanupRange.End)
$trailingChar.Delete()

# ---------------------------------------------------------------------------
# 2. Merge the two runs that together spell "Instructions for Mac OS Users:"
#    into a single run/text node.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("Instructions for Mac OS Users:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Instructions for Mac OS Users:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Remove the old "_GoBack" bookmark that used to live on the last
#    (otherwise empty) bullet item under the Mac OS instructions.
# ---------------------------------------------------------------------------

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()
